# Update cryptocurrency price/volume data per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.885.94"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "1.705.71"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'316.90"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "'0.4016"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").Value = "'0.4065"
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("D9").Value = "'1.483"
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("B10").Value = "BinanceUSD"
$ws.Range("C10").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D10").Value = "'1.003"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'53.61"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "'0.08813"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "'26.29"
$ws.Range("E13").Value = "  +4.35%  "
$ws.Range("D14").Value = "'7.476"
$ws.Range("E14").Value = "  -4.83%  "
$ws.Range("D15").Value = "'8.094"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").Value = "'0.00001355"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "1.723.05"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "'96.74"
$ws.Range("E18").Value = "  -4.09%  "
$ws.Range("D19").Value = "'0.07159"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'21.05"
$ws.Range("E20").Value = "  +3.12%  "
$ws.Range("D21").Value = "'7.248"
$ws.Range("E21").Value = "  -3.83%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").Value = "'14.37"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "24.895.77"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").Value = "'2.327"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").Value = "'2.896"
$ws.Range("E26").Value = "  -8.42%  "
$ws.Range("D27").Value = "'23.19"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").Value = "'6.286"
$ws.Range("E28").Value = "  +21.18%  "
$ws.Range("D29").Value = "'166.39"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "'146.01"
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("D31").Value = "'8.262"
$ws.Range("E31").Value = "  -11.58%  "
$ws.Range("D32").Value = "1.918.49"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").Value = "'2.230"
$ws.Range("E33").Value = "  +13.38%  "
$ws.Range("D34").Value = "'0.08877"
$ws.Range("E34").Value = "  -2.27%  "
$ws.Range("D35").Value = "'0.03204"
$ws.Range("E35").Value = "  +6.76%  "
$ws.Range("D36").Value = "'7.195"
$ws.Range("E36").Value = "  -8.64%  "
$ws.Range("D37").Value = "'1.025"
$ws.Range("E37").Value = "  -6.91%  "
$ws.Range("D38").Value = "'0.2853"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").Value = "'0.8469"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").Value = "'10.77"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").Value = "'0.09293"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'14.20"
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("D43").Value = "'1.473"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").Value = "'17.53"
$ws.Range("E44").Value = "  +5.47%  "
$ws.Range("D45").Value = "'2.727"
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").Value = "'0.7449"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "'4.246"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "'1.396"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("D49").Value = "'1.001"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "'142.13"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").Value = "'0.08358"
$ws.Range("E51").Value = "  +3.08%  "
